$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "22493944"
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "invictus10"
$ws.Range("E6").Value = "1234"
$ws.Range("F6").Value = "4321"
$ws.Range("G6").Value = "Acierto"
$ws.Range("H6").Value = "000"
$ws.Range("I6").Value = "0369"
$ws.Range("J6").Value = "NO ERROR"
$ws.Range("K6").Value = "bolp"
$ws.Range("L6").Value = "ACTIVO"
$ws.Range("M6").Value = "pruebasqa99"
$ws.Range("N6").Value = "jalzate@todo1.net"
$ws.Range("O6").Value = "Personal American Express"
$ws.Range("P6").Value = "*7806"
$ws.Range("Q6").Value = "Otro valor"
$ws.Range("R6").Value = "1000"
$ws.Range("S6").Value = "Pesos"

Write-Output "done"
